$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("D1").Value = "speciality"
$ws.Range("F1").Value = "state"
$ws.Range("G1").Value = "district"

# --- Row 2: Onkar Nagarkar ---
$ws.Range("B2").Value = "aid.techcryptors@gmail.com"
$ws.Range("D2").Value = "oncology"
$ws.Range("F2").Value = "Maharastra"
$ws.Range("G2").Value = "Mumbai"

# --- Row 3: Akanksha Indap ---
$ws.Range("B3").Value = "it.techcryptors@gmail.com"
$ws.Range("D3").Value = "sports-medicine"
$ws.Range("F3").Value = "Maharastra"
$ws.Range("G3").Value = "Raigad"

# --- Restore selected cell (matches the saved selection state in the diff) ---
$ws.Range("I8").Select()

Write-Host "Workbook updated"
